# Atualização automática: 2025-08-31 09:00:25
#
# This script rotates the detection data stored in rows 7-11 (the data that
# used to be on row 11 moves up to row 7, and rows 7,8,9,10 each shift down
# by one row), and refreshes the latest-detection fields on row 18 with a
# newer image/coords/confidence for the same fly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are stored as text in this sheet. Some of their values look
# numeric (e.g. "0.62" or "702,633,740,690"), so Excel's automatic type
# detection could silently turn them into real numbers. Prefixing with an
# apostrophe forces a text entry (exactly like typing it by hand), and
# resetting Style to "Normal" afterwards keeps the cell's formatting as it
# was (no explicit text-format / quote-prefix residue).
$textCols = @("A", "D", "E", "F", "I", "J")
# Columns that hold real numbers.
$numCols = @("G", "H")

function Set-CellText($range, [string]$val) {
    $range.Value = "'" + $val
    $range.Style = "Normal"
}

# --- Capture current values for the rotating block (rows 7-11) ---------
$cols = $textCols + $numCols

$snapshot = @{}
foreach ($r in 7..11) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# New order after the update: the data that was on row 11 becomes row 7,
# and rows 7,8,9,10 each shift down into 8,9,10,11.
$newOrder = @{
    7  = 11
    8  = 7
    9  = 8
    10 = 9
    11 = 10
}

foreach ($destRow in 7..11) {
    $srcRow = $newOrder[$destRow]
    $data = $snapshot[$srcRow]
    foreach ($c in $textCols) {
        Set-CellText $ws.Range("$c$destRow") $data[$c]
    }
    foreach ($c in $numCols) {
        $ws.Range("$c$destRow").Value = $data[$c]
    }
}

# --- Update row 18 with the newer detection image / coords / confidence -
Set-CellText $ws.Range("D18") "image_20250808221835_ppp0.jpg"
Set-CellText $ws.Range("I18") "1182,405,1231,455"
Set-CellText $ws.Range("J18") "0.76"

$wb.Save()
